# Change email to use send email library, remove unused email config.
#
# Constants sheet: remove the "BuyWaysDepartmentApproversEmailCredentials"
# row (row 19) -- the asset name used for logging into the old email
# credentials -- shifting the rows below it (BuyWaysLoginCredentials,
# BuyWaysLoginURLAssetName) up by one.
#
# Assets sheet: remove the "ExchangeServerURL" asset row (row 5), shifting
# the rows below it up by one.

$wb = $excel.ActiveWorkbook

$wsConstants = $wb.Worksheets.Item("Constants")
$wsConstants.Activate()
$wsConstants.Range("B30").Select() | Out-Null
$wsConstants.Rows.Item(19).Delete() | Out-Null

$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Activate()
$wsAssets.Range("B30").Select() | Out-Null
$wsAssets.Rows.Item(5).Delete() | Out-Null

$wsConstants.Activate() | Out-Null
